$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that are no longer present in the updated data (rows 6 and 7)
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# Update remaining data rows (2-5) with refreshed TPM-derived NATMI values
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Gm13306"
$ws.Range("C2").Value = "Ccr10"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5955496666666666
$ws.Range("H2").Value = 1.786649
$ws.Range("I2").Value = 0.5754499665999525
$ws.Range("J2").Value = 0.6273319775338025
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.1452775
$ws.Range("N2").Value = 0.290555
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.08651996669916666
$ws.Range("R2").Value = 0.519119800195
$ws.Range("S2").Value = 0.5754499665999525
$ws.Range("T2").Value = 0.6273319775338025
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Gm13306"
$ws.Range("C3").Value = "Ccr10"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01512033333333333
$ws.Range("H3").Value = 0.045361
$ws.Range("I3").Value = 0.01461002465226267
$ws.Range("J3").Value = 0.01592725030652961
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.1452775
$ws.Range("N3").Value = 0.290555
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.002196644225833333
$ws.Range("R3").Value = 0.013179865355
$ws.Range("S3").Value = 0.01461002465226267
$ws.Range("T3").Value = 0.01592725030652961
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Gm13306"
$ws.Range("C4").Value = "Ccr10"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.256774
$ws.Range("H4").Value = 0.513548
$ws.Range("I4").Value = 0.2481079211256428
$ws.Range("J4").Value = 0.1803180604576105
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.1452775
$ws.Range("N4").Value = 0.290555
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.037303484785
$ws.Range("R4").Value = 0.14921393914
$ws.Range("S4").Value = 0.2481079211256428
$ws.Range("T4").Value = 0.1803180604576105
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Gm13306"
$ws.Range("C5").Value = "Ccr10"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1674846666666666
$ws.Range("H5").Value = 0.502454
$ws.Range("I5").Value = 0.1618320876221421
$ws.Range("J5").Value = 0.1764227117020574
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.1452775
$ws.Range("N5").Value = 0.290555
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.02433175366166666
$ws.Range("R5").Value = 0.14599052197
$ws.Range("S5").Value = 0.1618320876221421
$ws.Range("T5").Value = 0.1764227117020574